{"js": "// Insert two new paragraphs right after the paragraph that contains\n// \"pip3 install dash-auth==1.3.2\":\n//   1) a paragraph with the text \"pip3 install dash-bootstrap-components\"\n//   2) a new empty paragraph\n//\n// This matches the diff, which shows the new paragraphs inserted\n// immediately after the \"dash-auth\" line and before the pre-existing\n// empty paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"pip3 install dash-auth==1.3.2\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === target) {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find paragraph: \" + target);\n}\n\n// Insert the \"dash-bootstrap-components\" paragraph directly after the\n// anchor paragraph.\nconst newTextPara = anchor.insertParagraph(\n  \"pip3 install dash-bootstrap-components\",\n  Word.InsertLocation.after\n);\n\n// Insert a new, empty paragraph directly after the one we just added.\nnewTextPara.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert two new paragraphs right after the paragraph that contains\n# \"pip3 install dash-auth==1.3.2\":\n#   1) a paragraph with the text \"pip3 install dash-bootstrap-components\"\n#   2) a new empty paragraph\n#\n# This matches the diff, which shows the new paragraphs inserted\n# immediately after the \"dash-auth\" line and before the pre-existing\n# empty paragraphs.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*pip3 install dash-auth==1.3.2*\") {\n        $target = $p\n    }\n}\n\nif ($target -ne $null) {\n    # Add an empty paragraph right after the anchor paragraph, then fill\n    # it in with the new pip install line.\n    $target.Range.InsertParagraphAfter()\n    $newTextPara = $target.Next()\n    $newTextPara.Range.InsertBefore(\"pip3 install dash-bootstrap-components\")\n\n    # Add a second, empty paragraph right after the one we just filled in.\n    $newTextPara.Range.InsertParagraphAfter()\n}\n"}
